$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove workbook protection element (matches removal of <workbookProtection/> in workbook.xml)
try { $wb.Unprotect() } catch { }

$ws.Activate()

# New log rows appended by the "Add files via upload" commit
$rows = @(
    @(7,  "2024-09-29 13:36:49", "tickets",     879),
    @(8,  "2024-09-29 14:24:21", "evaluations", 26346),
    @(9,  "2024-09-29 14:44:47", "users",       173),
    @(10, "2024-09-29 15:27:33", "units",       97),
    @(11, "2024-09-29 21:52:32", "items",       196834),
    @(12, "2024-09-29 22:13:26", "items",       0),
    @(13, "2024-09-29 22:14:15", "categories",  21308),
    @(14, "2024-09-29 22:14:18", "actions",     13),
    @(15, "2024-09-29 22:14:19", "tickets",     1),
    @(16, "2024-09-29 22:15:00", "evaluations", 22),
    @(17, "2024-09-29 22:15:01", "users",       0),
    @(18, "2024-09-29 22:15:01", "units",       0),
    @(19, "2024-09-29 22:15:01", "departments", 7),
    @(20, "2024-09-29 22:15:02", "user_types",  26)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
}

# Reflect the new active selection / scroll position recorded in the sheetView
$ws.Range("O10").Select() | Out-Null
